# Translate the ContosoLearn Competitor SWOT document from English to Spanish.
#
# Each SWOT paragraph has the shape "<Bold Label:><space><body sentences>".
# We replace the bold label and the body text of each paragraph by
# addressing sub-ranges of the paragraph's Range directly (rather than
# Find/Replace) so that:
#   - run formatting (bold label vs. regular body) is preserved exactly,
#   - label text that repeats across paragraphs (e.g. "Strengths:") is only
#     ever touched within the paragraph currently being processed,
#   - no "smart quote" autocorrection mangles the literal '"' character that
#     appears at the end of the document.

$d = $word.ActiveDocument

# Replaces the bold "<Label>:" run and the regular body run of a two-run
# SWOT paragraph in a single step, using the length of the *original*
# label text to find where the body run starts (so stale offsets never
# leak in), then writing both sub-ranges back-to-front so that changing
# the label's length never invalidates the body range we already computed.
function Set-LabelAndBody($paraIndex, $oldLabelLength, $newLabel, $newBody) {
    $full = $d.Paragraphs($paraIndex).Range
    $bodyRange = $d.Range($full.Start + $oldLabelLength, $full.End)
    $bodyRange.Text = $newBody

    $labelRange = $d.Range($full.Start, $full.Start + $oldLabelLength)
    $labelRange.Text = $newLabel
}

# --- Title (paragraph 1) ---
$d.Paragraphs(1).Range.Find.Execute("ContosoLearn Competitor SWOT", $true, $false, $false, $false, $false, $true, 0, $false, "DAFO de competidores ContosoLearn", 2)

# --- Fabrikam Learning: (paragraph 2 heading stays the same) ---

# Paragraph 3: Strengths
Set-LabelAndBody 3 10 "Fortalezas:" ' Fabrikam Learning proporciona un conjunto completo de herramientas de análisis e informes. Garantiza la supervisión continua de las actividades de enseñanza y aprendizaje, así como la identificación de áreas problemáticas que deben abordarse.'

# Paragraph 4: Weaknesses
Set-LabelAndBody 4 11 "Puntos débiles:" ' aunque Fabrikam Learning tiene funcionalidades de informes sólidas, puede ser abrumador para algunos usuarios debido a su naturaleza completa.'

# Paragraph 5: Opportunities
Set-LabelAndBody 5 14 "Oportunidades:" ' hay una creciente demanda de experiencias de aprendizaje personalizadas y recomendaciones controladas por datos. Fabrikam Learning puede aprovechar sus sólidas herramientas de análisis e informes para satisfacer esta demanda.'

# Paragraph 6: Threats
Set-LabelAndBody 6 8 "Amenazas:" ' el mercado de eLearning es altamente competitivo con muchos jugadores que ofrecen características similares. Fabrikam Learning debe innovar continuamente para mantenerse a la vanguardia.'

# --- AdatumLearn: (paragraph 7 heading stays the same) ---

# Paragraph 8: Strengths
Set-LabelAndBody 8 10 "Fortalezas:" ' AdatumLearn ofrece cursos sobre técnicas de análisis de negocios como MOST y DAFO. Esto muestra su compromiso de proporcionar contenido valioso a sus usuarios.'

# Paragraph 9: Weaknesses
Set-LabelAndBody 9 11 "Puntos débiles:" ' la información proporcionada en sus cursos es una compilación de información generada por terceros. Esto podría no ser tan valioso como el contenido original.'

# Paragraph 10: Opportunities
Set-LabelAndBody 10 14 "Oportunidades:" ' AdatumLearn puede crear contenido más original para proporcionar un valor único a sus usuarios. También puede ampliar sus ofertas de cursos para tratar más temas.'

# Paragraph 11: Threats
Set-LabelAndBody 11 8 "Amenazas:" ' al igual que Fabrikam Learning, AdatumLearn también se enfrenta a una competencia rígida en el mercado de eLearning. Necesita mejorar continuamente su oferta para mantenerse competitivo".'
